$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "64.964.51"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.65%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.542.55"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.32%  "

$ws.Range("E4").Value = "  -0.10%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "581.35"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.41%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "152.94"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +4.15%  "

$ws.Range("E7").Value = "  +0.05%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.539"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.25%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.542.12"
$c.Style = "Normal"

$ws.Range("E10").Value = "  +2.09%  "

$ws.Range("E11").Value = "  -1.67%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.29"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.69%  "

$ws.Range("E13").Value = "  +0.64%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "29.29"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.06%  "

$ws.Range("E15").Value = "  +2.86%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.997.48"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.06%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "64.294.26"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.71%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.540.31"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.85%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "8.00"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.84%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.02"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.76%  "

$ws.Range("E21").Value = "  +3.69%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "329.32"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.05%  "

$ws.Range("E23").Value = "  +3.91%  "

$ws.Range("E24").Value = "  +0.03%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "10.25"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.68%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "65.87"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.76%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "645.21"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +1.20%  "

$ws.Range("E28").Value = "  +7.74%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.666.03"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.28%  "

$ws.Range("E30").Value = "  +4.68%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.16%  "

$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("E33").Value = "  +2.72%  "

$ws.Range("E34").Value = "  +2.24%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.04%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.57"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.19%  "

$ws.Range("E37").Value = "  +2.78%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "5.61"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +5.68%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "155.69"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.60%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.85"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.02%  "

$ws.Range("E41").Value = "  +1.47%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "18.95"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.87%  "

$ws.Range("E43").Value = "  +5.57%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "162.07"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +6.13%  "

$ws.Range("E45").Value = "  +0.05%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0₆0301"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("E47").Value = "  +2.59%  "

$ws.Range("E48").Value = "  +2.86%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "21.51"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +6.21%  "

$ws.Range("E50").Value = "  +2.90%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0519"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.33%  "
